$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.010.35'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.593.04'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.83'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.479'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.247'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0612'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.03'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  +2.24%  '
$ws.Range('D12').Value = '1.813.83'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '1.584.16'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.513'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '26.071.31'
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.48'
$ws.Range('E17').Value = '  +2.91%  '
$ws.Range('D18').Value = '0.0₃0724'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '202.26'
$ws.Range('E20').Value = '  +5.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.25'
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.26'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.99'
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  +15.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.74'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.121'
$ws.Range('E27').Value = '  -7.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.20'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.50'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.16'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0475'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.11'
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.90'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = '1.129.56'
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0162'
$ws.Range('E37').Value = '  +7.96%  '
$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.793'
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.31'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.494'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.782'
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.15'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '1.729.49'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.05'
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.49'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.71'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₇0972'
$ws.Range('E49').Value = '  -13.15%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.407'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.08%  '
